$d = $word.ActiveDocument

$replacements = @(
    @{old="11×12=132"; new="29×23=667"},
    @{old="30×26=780"; new="71×31=2201"},
    @{old="41×35=1435"; new="81×70=5670"},
    @{old="62×56=3472"; new="91×50=4550"},
    @{old="20×46=920"; new="62×21=1302"},
    @{old="35×34=1190"; new="52×62=3224"},
    @{old="63×18=1134"; new="93×16=1488"},
    @{old="69×26=1794"; new="20×32=640"},
    @{old="57×67=3819"; new="79×66=5214"},
    @{old="37×69=2553"; new="26×87=2262"},
    @{old="68×30=2040"; new="81×42=3402"},
    @{old="72×68=4896"; new="17×80=1360"},
    @{old="47×96=4512"; new="68×39=2652"},
    @{old="33×21=693"; new="78×90=7020"},
    @{old="42×94=3948"; new="58×21=1218"},
    @{old="57×64=3648"; new="24×13=312"},
    @{old="95×78=7410"; new="26×65=1690"},
    @{old="65×54=3510"; new="13×20=260"},
    @{old="14×15=210"; new="39×75=2925"},
    @{old="88×95=8360"; new="48×75=3600"},
    @{old="42×53=2226"; new="69×23=1587"},
    @{old="98×67=6566"; new="15×74=1110"},
    @{old="74×84=6216"; new="21×84=1764"},
    @{old="34×63=2142"; new="25×44=1100"},
    @{old="32×50=1600"; new="68×76=5168"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
